# calorimetry : input and output consistency : done
#
# The free components (H, L, M) were mistakenly carried into the
# per-species *result* tables alongside the actual fitted complexes
# (HL, H2L, ML, HML). Drop those rows so inputs/outputs line up, and
# correct the sign on the recomputed reaction enthalpies for ML / HML.

$wb = $excel.ActiveWorkbook

# --- constants_evaluated: drop H, L, M rows (2:4), keep HL/H2L/ML/HML ---
$wsConstants = $wb.Worksheets.Item("constants_evaluated")
$wsConstants.Range("A2:A4").EntireRow.Delete()

# --- input_enthalpies: drop H, L, M rows (2:4), keep HL/H2L ---
$wsInputEnthalpies = $wb.Worksheets.Item("input_enthalpies")
$wsInputEnthalpies.Range("A2:A4").EntireRow.Delete()

# --- enthalpies_calculated: drop H, L, M rows (2:4), keep HL/H2L/ML/HML,
#     and fix the sign of the recomputed reaction enthalpy for ML/HML ---
$wsEnthalpiesCalc = $wb.Worksheets.Item("enthalpies_calculated")
$wsEnthalpiesCalc.Range("A2:A4").EntireRow.Delete()
$wsEnthalpiesCalc.Range("B4").Value = -9.49363111862915
$wsEnthalpiesCalc.Range("B5").Value = -18.9066684805445
